# Scheduled runner update: refresh computed market-price / profit columns
# (H, I, J, K, L, M, N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets,
# leaving the Leve metadata columns (A-G) untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1066.4717
$ws.Range("J17").Value = 874.2308
$ws.Range("L17").Value = 2622.6924
$ws.Range("N17").Value = -2958.6924

$ws.Range("H62").Value = 2610.375
$ws.Range("I62").Value = 1779.4
$ws.Range("K62").Value = 1779.4
$ws.Range("M62").Value = -1155.4

$ws.Range("H65").Value = 2610.375
$ws.Range("I65").Value = 1779.4
$ws.Range("K65").Value = 8897
$ws.Range("M65").Value = -5777

$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996

$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984

$ws.Range("H132").Value = 984.65
$ws.Range("I132").Value = 989.38464
$ws.Range("K132").Value = 2968.15392
$ws.Range("M132").Value = -438.1539199999997

$ws.Range("H138").Value = 1790.1316
$ws.Range("I138").Value = 1787.7646
$ws.Range("J138").Value = 1792.0476
$ws.Range("K138").Value = 5363.293799999999
$ws.Range("L138").Value = 5376.142800000001
$ws.Range("M138").Value = -223.2937999999995
$ws.Range("N138").Value = -15656.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 758.1667
$ws.Range("I2").Value = 499.875
$ws.Range("K2").Value = 499.875
$ws.Range("M2").Value = -386.875

$ws.Range("H61").Value = 2869.2354
$ws.Range("I61").Value = 1530.8572
$ws.Range("K61").Value = 1530.8572
$ws.Range("M61").Value = -1318.8572

$ws.Range("H116").Value = 758.1667
$ws.Range("I116").Value = 499.875
$ws.Range("K116").Value = 499.875
$ws.Range("M116").Value = 1794.125

$ws.Range("H136").Value = 2869.2354
$ws.Range("I136").Value = 1530.8572
$ws.Range("K136").Value = 4592.571599999999
$ws.Range("M136").Value = -2042.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 758.1667
$ws.Range("I3").Value = 499.875
$ws.Range("K3").Value = 499.875
$ws.Range("M3").Value = -385.875

$ws.Range("H105").Value = 1826.4
$ws.Range("I105").Value = 1951.125
$ws.Range("J105").Value = 1327.5
$ws.Range("K105").Value = 1951.125
$ws.Range("L105").Value = 1327.5
$ws.Range("M105").Value = -204.125
$ws.Range("N105").Value = -4821.5

$ws.Range("H107").Value = 1641.75
$ws.Range("I107").Value = 1641.75
$ws.Range("K107").Value = 1641.75
$ws.Range("M107").Value = 278.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1175
$ws.Range("J22").Value = 1433.3334
$ws.Range("L22").Value = 1433.3334
$ws.Range("N22").Value = -2133.3334

$ws.Range("H31").Value = 3389.7896
$ws.Range("I31").Value = 1482.6666
$ws.Range("K31").Value = 1482.6666
$ws.Range("M31").Value = -1187.6666

$ws.Range("H34").Value = 3389.7896
$ws.Range("I34").Value = 1482.6666
$ws.Range("K34").Value = 1482.6666
$ws.Range("M34").Value = -1280.6666

$ws.Range("H58").Value = 1196.1333
$ws.Range("I58").Value = 988.1177
$ws.Range("J58").Value = 1468.1538
$ws.Range("K58").Value = 988.1177
$ws.Range("L58").Value = 1468.1538
$ws.Range("M58").Value = -785.1177
$ws.Range("N58").Value = -1874.1538

$ws.Range("H105").Value = 1505.75
$ws.Range("I105").Value = 1006.5714
$ws.Range("K105").Value = 1006.5714
$ws.Range("M105").Value = 740.4286

$ws.Range("H136").Value = 1196.1333
$ws.Range("I136").Value = 988.1177
$ws.Range("J136").Value = 1468.1538
$ws.Range("K136").Value = 2964.3531
$ws.Range("L136").Value = 4404.4614
$ws.Range("M136").Value = -414.3531000000003
$ws.Range("N136").Value = -9504.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 654.9091
$ws.Range("I5").Value = 588.7778
$ws.Range("K5").Value = 1766.3334
$ws.Range("M5").Value = -1654.3334

$ws.Range("H107").Value = 727.8125
$ws.Range("I107").Value = 295.5
$ws.Range("K107").Value = 886.5
$ws.Range("M107").Value = 1033.5

$ws.Range("H131").Value = 787.3
$ws.Range("I131").Value = 432.33334
$ws.Range("J131").Value = 798.2782999999999
$ws.Range("K131").Value = 1297.00002
$ws.Range("L131").Value = 2394.8349
$ws.Range("M131").Value = 3742.99998
$ws.Range("N131").Value = -12474.8349

$ws.Range("H135").Value = 654.9091
$ws.Range("I135").Value = 588.7778
$ws.Range("K135").Value = 5299.000199999999
$ws.Range("M135").Value = -2764.000199999999

$ws.Range("H139").Value = 10891.8
$ws.Range("I139").Value = 10891.8
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 32675.4
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -27535.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 8512478
$ws.Range("J11").Value = 6713336
$ws.Range("L11").Value = 6713336
$ws.Range("N11").Value = -6713614

$ws.Range("H113").Value = 1922.2
$ws.Range("I113").Value = 1605.5
$ws.Range("K113").Value = 1605.5
$ws.Range("M113").Value = 564.5

$ws.Range("H122").Value = 2261.7646
$ws.Range("J122").Value = 2681.3333
$ws.Range("L122").Value = 8043.999899999999
$ws.Range("N122").Value = -12943.9999

$ws.Range("H132").Value = 5192.24
$ws.Range("I132").Value = 4220.1665
$ws.Range("K132").Value = 12660.4995
$ws.Range("M132").Value = -10130.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 8499.333000000001
$ws.Range("J14").Value = 8499.333000000001
$ws.Range("L14").Value = 8499.333000000001
$ws.Range("N14").Value = -8843.333000000001

$ws.Range("H46").Value = 1583.3334
$ws.Range("J46").Value = 1583.3334
$ws.Range("L46").Value = 1583.3334
$ws.Range("N46").Value = -1959.3334

$ws.Range("H61").Value = 4625
$ws.Range("I61").Value = 4250
$ws.Range("K61").Value = 4250
$ws.Range("M61").Value = -4048

$ws.Range("H113").Value = 4625
$ws.Range("I113").Value = 4250
$ws.Range("K113").Value = 4250
$ws.Range("M113").Value = -2080

$ws.Range("H136").Value = 2855.3215
$ws.Range("I136").Value = 1678.9231
$ws.Range("J136").Value = 3874.8667
$ws.Range("K136").Value = 5036.7693
$ws.Range("L136").Value = 11624.6001
$ws.Range("M136").Value = -2486.7693
$ws.Range("N136").Value = -16724.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1180.6666
$ws.Range("I113").Value = 1134.6666
$ws.Range("K113").Value = 3403.9998
$ws.Range("M113").Value = -1233.9998

$ws.Range("H126").Value = 5882.44
$ws.Range("I126").Value = 5223.0625
$ws.Range("K126").Value = 15669.1875
$ws.Range("M126").Value = -13199.1875

$ws.Range("H132").Value = 3590.3333
$ws.Range("I132").Value = 3258.9412
$ws.Range("K132").Value = 9776.8236
$ws.Range("M132").Value = -7246.8236

$ws.Range("H136").Value = 3956.5715
$ws.Range("I136").Value = 3898.2
$ws.Range("K136").Value = 11694.6
$ws.Range("M136").Value = -9144.599999999999
